# Update for new taxotools:
# Clear the family (A), genus (C), species (D) and author (E) columns
# for every data row back to the sheet's "NA" placeholder, since this
# information is now produced by the refreshed taxotools pipeline
# downstream (specificEpithet/infraspecificEpithet/taxonRank/
# scientificNameAuthorship/canonicalName columns already carry it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "NA"   # family
    $ws.Cells.Item($r, 3).Value = "NA"   # genus
    $ws.Cells.Item($r, 4).Value = "NA"   # species
    $ws.Cells.Item($r, 5).Value = "NA"   # author
}

# Reset the view: scroll back to the top-left and select A2, matching
# the reviewed workbook's saved cursor position.
$ws.Range("A2").Select() | Out-Null
